$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$d3 = @'
Background
id="Par1">Early in the pandemic, inadequate SARS-CoV-2 testing limited understanding of transmission.

 Chief among barriers to large-scale testing was unknown feasibility, particularly in non-urban areas.

 Our objective was to report methods of high-volume, comprehensive SARS-CoV-2 testing, offering one model to augment disease surveillance in a rural community.


Methods
id="Par2">A community-university partnership created an operational site used to test most residents of Bolinas, California regardless of symptoms in 4 days (April 20th – April 23rd, 2020).

 Prior to testing, key preparatory elements included community mobilization, pre-registration, volunteer recruitment, and data management.

 On day of testing, participants were directed to a testing lane after site entry.

 An administrator viewed the lane-specific queue and pre-prepared test kits, linked to participants’ records.

 Medical personnel performed sample collection, which included finger prick with blood collection to run laboratory-based antibody testing and respiratory specimen collection for polymerase chain reaction (PCR).


Results
id="Par3">Using this 4-lane model, 1,840 participants were tested in 4 days.

 A median of 57 participants (IQR 47–67) were tested hourly.

 The fewest participants were tested on day 1 (n = 338 participants), an intentionally lower volume day, increasing to n = 571 participants on day 4. The number of testing teams was also increased to two per lane to allow simultaneous testing of multiple participants on days 2–4. Consistent staffing on all days helped optimize proficiency, and strong community partnership was essential from planning through execution.


Conclusions
id="Par4">High-volume ascertainment of SARS-CoV-2 prevalence by PCR and antibody testing was feasible when conducted in a community-led, drive-through model in a non-urban area.


Supplementary Information
The online version contains supplementary material available at 10.1186/s13690-021-00647-8.

'@

$e3 = @'
[Ayesha%Appa%Ayesha.appa@ucsf.edu%1,  Gabriel%Chamie%NULL%1,  Aenor%Sawyer%NULL%1,  Kimberly%Baltzell%NULL%1,  Kathryn%Dippell%NULL%1,  Salu%Ribeiro%NULL%1,  Elias%Duarte%NULL%1,  Joanna%Vinden%NULL%1,  CLIAHUB%Consortium%NULL%1,  Jonathan%Kramer-Feldman%NULL%1,  Shahryar%Rahdari%NULL%1,  Doug%MacIntosh%NULL%1,  Katherine%Nicholson%NULL%1,  Jonathan%Im%NULL%1,  Diane%Havlir%NULL%1,  Bryan%Greenhouse%NULL%0]
'@

$d15 = @'
id="Par1">In Taiwan, high-risk patients have been identified and tested for preventing community spread of COVID-19. Most sample collection was performed in emergency departments (EDs).

 Traditional sample collection requires substantial personal protective equipment (PPE), healthcare professionals, sanitation workers, and isolation space.

 To solve this problem, we established a multifunctional sample collection station (MSCS) for COVID-19 testing in front of our ED.

 The station is composed of a thick and clear acrylic board (2 cm), which completely separates the patient and medical personnel.

 Three pairs of gloves (length, 45 cm) are attached and fixed on the outside wall of the MSCS.

 The gloves are used to conduct sampling of throat/nasal swabs, sputum, and blood from patients.

 The gap between the board and the building is only 0.2 cm (sealed with silicone sealant).

 ED personnel communicate with patients using a small two-way broadcast system.

 Medical waste is put in specific trashcans installed in the table outside the MSCS.

 With full physical protection, the personnel conducting the sampling procedure need to wear only their N95 mask and gloves.

 After we activated the station, our PPE, sampling time, and sanitization resources were considerably conserved during the 4-week observation period.

 The MSCS obviously saved time and PPE.

 It elevated the efficiency and capacity of the ED for handling potential community infections of COVID-19.
'@

$e15 = @'
[Po-Ting%Lin%NULL%1,  Ting-Yuan%Ni%NULL%1,  Tren-Yi%Chen%NULL%1,  Chih-Pei%Su%NULL%1,  Hsiao-Fen%Sun%NULL%1,  Mu-Kuan%Chen%NULL%1,  Chu-Chung%Chou%NULL%1,  Po-Yu%Wang%NULL%1,  Yan-Ren%Lin%H6213.lac@gmail.com%1]
'@

$ws.Range("D3").Value = $d3
$ws.Range("E3").Value = $e3
$ws.Range("D15").Value = $d15
$ws.Range("E15").Value = $e15
